# TestNG Listeners are added
# Update the OrangeHRM test-data sheet row 5 (TC No "5") sample
# values in the Message / Middle Name / LastName columns (F, H, I)
# from the old placeholder "cfmYGy" to the new placeholder "zxZYll".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrangeHRM")

$ws.Range("F6").Value = "zxZYll"
$ws.Range("H6").Value = "zxZYll"
$ws.Range("I6").Value = "zxZYll"
